# feat: Separar gastos de obra do budget mensal
# Applies the numeric/text updates described in the commit:
#  - refresh "Atualizado"/sync timestamps
#  - Dashboard: Gastos Variaveis (C7/D7) and Obra (C9) KPI rows
#  - Dashboard "Saude" category row (C16/D16)
#  - Mensal sheet "Transporte" monthly value (B6)
#  - Categorias sheet "Transporte" row (C6/D6/E6)
#  - Dados sheet sync_timestamp (B3) and "saude" gasto_jan (D10)

$wb = $excel.ActiveWorkbook

# ---- Dashboard ----
$dash = $wb.Worksheets.Item("Dashboard")

# Updated timestamp banner
$dash.Range("A2").Value2 = "Atualizado: 31/12/2025 11:09"

# Resumo do Mes: Gastos Variaveis row (real + variacao%)
$dash.Range("C7").Value2 = 9265.58
$dash.Range("D7").Value = "'-47%"

# Resumo do Mes: Obra row (real)
$dash.Range("C9").Value2 = 36234.42

# Gastos por Categoria: Saude row (real + %)
$dash.Range("C16").Value2 = 1730.24
$dash.Range("D16").Value2 = 346

# ---- Mensal ----
$mensal = $wb.Worksheets.Item("Mensal")
$mensal.Range("B6").Value2 = 1730.24

# ---- Categorias ----
$categorias = $wb.Worksheets.Item("Categorias")
$categorias.Range("C6").Value2 = 1730.24
$categorias.Range("D6").Value2 = -1230.24
$categorias.Range("E6").Value2 = 3.46048

# ---- Dados ----
$dados = $wb.Worksheets.Item("Dados")
$dados.Range("B3").Value2 = "2025-12-31T11:09:01.247836"
$dados.Range("D10").Value2 = 1730.24
